# Add new columns I0 (col I) and IF (col J) to the activity log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) -------------------------------------------------
# Copy the formatting of the existing last header cell (H1, style index 1:
# bold font + border + center/top alignment) onto the two new header cells
# so they match the rest of the header row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data cells (rows 2-76) ------------------------------------------------
$iVals = @(9,9,9,9,9,8,9,9,9,9,9,9,8,9,9,9,9,9,9,9,9,9,10,9,9,8,8,9,9,8,8,8,8,8,7,7,7,8,7,9,8,8,8,8,7,9,8,9,6,9,8,9,9,9,9,9,9,8,9,9,9,9,8,9,9,9,9,6,8,8,7,7,8,6,7)
$jVals = @(9,9,9,9,9,8,9,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,8,8,9,9,9,9,8,8,8,8,8,8,8,8,9,8,8,8,9,8,9,8,9,7,9,8,9,9,9,9,9,9,8,9,9,9,9,8,9,9,9,9,7,8,8,7,7,8,6,7)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value  = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}

Write-Output "I0 and IF columns added"
